# Update the sample/test data on the "herzo_student" worksheet.
# Row 2 and Row 3 hold raw given-name / surname values (columns A and B);
# columns D, E and F are derived via formulas and recalculate automatically.
# Column G holds the generated password, which is also replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: first test user
$ws.Range("A2").Value = "test1gn"
$ws.Range("B2").Value = "test1sn"

# Row 3: second test user
$ws.Range("A3").Value = "test2gn"
$ws.Range("B3").Value = "test2sn"

# New shared password for both generated accounts
$ws.Range("G2").Value = "Q121ghhjg!a"
$ws.Range("G3").Value = "Q121ghhjg!a"

# Recalculate so the formula-driven columns (D, E, F) pick up the new values
$excel.Calculate()

# Restore the on-screen selection to the last edited cell
$null = $ws.Range("G3").Select()
